# "Cube with textures added" - mark Milestone II items complete for several
# rubric rows, flip the Milestone II Complete(X) flag on, and update the
# selected cell to reflect where the editor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that simply record "II" as the milestone this feature was completed on.
$milestoneIIRows = 20, 29, 30, 31, 46, 47, 52, 55
foreach ($r in $milestoneIIRows) {
    $ws.Range("E$r").Value = "II"
}

# Row 6: mark milestone "II" and flag it done with an "X" in F6.
$ws.Range("E6").Value = "II"
$ws.Range("F6").Value = "X"

# Milestone II Complete(X) flags for the citation rows.
$ws.Range("D73").Value = "X"
$ws.Range("D74").Value = "X"

# Restore the selection to where the editor left off.
$ws.Range("F55").Select()

$wb.Save()
